# Added Test Data for Croatia Market
#
# - Duplicates the "Slovakia" sheet (same layout/styles) to the end of the
#   workbook, renames the copy to "Croatia", and fills in the
#   market name / ticket reference for Croatia.
# - Leaves the new "Croatia" tab selected/active (last tab), which also
#   clears the previously-active "Turkey" tab's tabSelected flag.
# - Marks the "Slovakia" sheet's selection as the whole sheet
#   (mirrors a `Ctrl+A` / select-all that was left behind on that tab).

$wb = $excel.ActiveWorkbook

# 1. Update Slovakia's remembered selection to "select all" before we branch
#    off its copy, so the copy does not inherit this.
$slovakia = $wb.Worksheets.Item("Slovakia")
$slovakia.Range("A1:XFD1048576").Select() | Out-Null

# 2. Duplicate Slovakia -> new sheet at the very end of the workbook
#    (placed immediately after the current last tab, "Turkey").
$turkey = $wb.Worksheets.Item("Turkey")
$slovakia.Copy([System.Reflection.Missing]::Value, $turkey)
$croatia = $wb.Worksheets.Item($wb.Worksheets.Count)
$croatia.Name = "Croatia"

# 3. Fill in the Croatia-specific data.
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2473"

# 4. Leave the selection on the new sheet at B4, and make it the active tab
#    (this Select() happens last, so it becomes the workbook's active sheet).
$croatia.Range("B4").Select() | Out-Null
